$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 7978
$ws.Range("I51").Value = 2995
$ws.Range("K51").Value = 2995
$ws.Range("M51").Value = -2511
$ws.Range("H76").Value = 4926.636
$ws.Range("J76").Value = 4899.625
$ws.Range("L76").Value = 4899.625
$ws.Range("N76").Value = -5529.625
$ws.Range("H79").Value = 4926.636
$ws.Range("J79").Value = 4899.625
$ws.Range("L79").Value = 4899.625
$ws.Range("N79").Value = -7083.625
$ws.Range("H101").Value = 3005.5
$ws.Range("I101").Value = 1341.6666
$ws.Range("K101").Value = 4024.9998
$ws.Range("M101").Value = -2402.9998
$ws.Range("H138").Value = 3897.0193
$ws.Range("I138").Value = 4032.75
$ws.Range("J138").Value = 3885.7083
$ws.Range("K138").Value = 12098.25
$ws.Range("L138").Value = 11657.1249
$ws.Range("M138").Value = -6958.25
$ws.Range("N138").Value = -21937.1249
$ws.Range("H141").Value = 7733.3
$ws.Range("I141").Value = 6451.7646
$ws.Range("K141").Value = 19355.2938
$ws.Range("M141").Value = -14175.2938

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3467.0254
$ws.Range("I32").Value = 1780.8939
$ws.Range("K32").Value = 1780.8939
$ws.Range("M32").Value = -1493.8939
$ws.Range("H88").Value = 125749.5
$ws.Range("J88").Value = 125749.5
$ws.Range("L88").Value = 125749.5
$ws.Range("N88").Value = -126561.5
$ws.Range("H91").Value = 125749.5
$ws.Range("J91").Value = 125749.5
$ws.Range("L91").Value = 125749.5
$ws.Range("N91").Value = -128557.5
$ws.Range("H110").Value = 1574785.1
$ws.Range("I110").Value = 2553402.2
$ws.Range("J110").Value = 8997.6
$ws.Range("K110").Value = 2553402.2
$ws.Range("L110").Value = 8997.6
$ws.Range("M110").Value = -2551357.2
$ws.Range("N110").Value = -13087.6

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2205.7368
$ws.Range("I20").Value = 2598.5715
$ws.Range("J20").Value = 1105.8
$ws.Range("K20").Value = 2598.5715
$ws.Range("L20").Value = 1105.8
$ws.Range("M20").Value = -2351.5715
$ws.Range("N20").Value = -1599.8
$ws.Range("H22").Value = 502.44446
$ws.Range("I22").Value = 502.44446
$ws.Range("K22").Value = 502.44446
$ws.Range("M22").Value = -329.44446
$ws.Range("H86").Value = 1165.8334
$ws.Range("I86").Value = 1139.9412
$ws.Range("K86").Value = 1139.9412
$ws.Range("M86").Value = -16.94119999999998
$ws.Range("H89").Value = 1165.8334
$ws.Range("I89").Value = 1139.9412
$ws.Range("K89").Value = 5699.706
$ws.Range("M89").Value = -83.70600000000013
$ws.Range("H105").Value = 2925.64
$ws.Range("J105").Value = 3227.6924
$ws.Range("L105").Value = 3227.6924
$ws.Range("N105").Value = -6721.6924
$ws.Range("H107").Value = 1985.2941
$ws.Range("I107").Value = 1194.5927
$ws.Range("J107").Value = 5035.143
$ws.Range("K107").Value = 1194.5927
$ws.Range("L107").Value = 5035.143
$ws.Range("M107").Value = 725.4073000000001
$ws.Range("N107").Value = -8875.143

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 607
$ws.Range("I22").Value = 321.16666
$ws.Range("J22").Value = 950
$ws.Range("K22").Value = 321.16666
$ws.Range("L22").Value = 950
$ws.Range("M22").Value = 28.83334000000002
$ws.Range("N22").Value = -1650
$ws.Range("H94").Value = 2895
$ws.Range("I94").Value = 2435.6
$ws.Range("K94").Value = 2435.6
$ws.Range("M94").Value = -1984.6
$ws.Range("H99").Value = 14323.615
$ws.Range("I99").Value = 23118.5
$ws.Range("J99").Value = 6785.143
$ws.Range("K99").Value = 23118.5
$ws.Range("L99").Value = 6785.143
$ws.Range("M99").Value = -21620.5
$ws.Range("N99").Value = -9781.143
$ws.Range("H105").Value = 910390.6
$ws.Range("I105").Value = 1338075
$ws.Range("K105").Value = 1338075
$ws.Range("M105").Value = -1336328
$ws.Range("H126").Value = 14323.615
$ws.Range("I126").Value = 23118.5
$ws.Range("J126").Value = 6785.143
$ws.Range("K126").Value = 69355.5
$ws.Range("L126").Value = 20355.429
$ws.Range("M126").Value = -66885.5
$ws.Range("N126").Value = -25295.429

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 873.24243
$ws.Range("I107").Value = 1168.3
$ws.Range("K107").Value = 3504.9
$ws.Range("M107").Value = -1584.9
$ws.Range("H126").Value = 9407.5
$ws.Range("I126").Value = 9407.5
$ws.Range("K126").Value = 28222.5
$ws.Range("M126").Value = -23282.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5296728.5
$ws.Range("I70").Value = 11910014
$ws.Range("K70").Value = 11910014
$ws.Range("M70").Value = -11909744
$ws.Range("H73").Value = 5296728.5
$ws.Range("I73").Value = 11910014
$ws.Range("K73").Value = 11910014
$ws.Range("M73").Value = -11909078
$ws.Range("H102").Value = 22734210
$ws.Range("I102").Value = 29419956
$ws.Range("J102").Value = 2677
$ws.Range("K102").Value = 29419956
$ws.Range("L102").Value = 2677
$ws.Range("M102").Value = -29418334
$ws.Range("N102").Value = -5921
$ws.Range("H113").Value = 3037.6
$ws.Range("I113").Value = 3037.6
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3037.6
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -867.5999999999999
$ws.Range("H122").Value = 2989.8333
$ws.Range("I122").Value = 2249.5
$ws.Range("K122").Value = 6748.5
$ws.Range("M122").Value = -4298.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1021
$ws.Range("I9").Value = 1021
$ws.Range("K9").Value = 1021
$ws.Range("M9").Value = -797
$ws.Range("H55").Value = 233.14285
$ws.Range("J55").Value = 223.25
$ws.Range("L55").Value = 223.25
$ws.Range("N55").Value = -569.25
$ws.Range("H82").Value = 4466068.5
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 4466068.5
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H100").Value = 2499.75
$ws.Range("H132").Value = 3991.9194
$ws.Range("I132").Value = 2932.422
$ws.Range("K132").Value = 8797.266
$ws.Range("M132").Value = -6267.266
$ws.Range("H136").Value = 2546.3462
$ws.Range("I136").Value = 2076.2856
$ws.Range("K136").Value = 6228.8568
$ws.Range("M136").Value = -3678.8568

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3479787.8
$ws.Range("I81").Value = 2612541
$ws.Range("J81").Value = 5214281
$ws.Range("K81").Value = 5225082
$ws.Range("L81").Value = 10428562
$ws.Range("M81").Value = -5224021
$ws.Range("N81").Value = -10430684
$ws.Range("H84").Value = 3479787.8
$ws.Range("I84").Value = 2612541
$ws.Range("J84").Value = 5214281
$ws.Range("K84").Value = 26125410
$ws.Range("L84").Value = 52142810
$ws.Range("M84").Value = -26120106
$ws.Range("N84").Value = -52153418
$ws.Range("H107").Value = 2754.2942
$ws.Range("I107").Value = 3269.5833
$ws.Range("J107").Value = 1517.6
$ws.Range("K107").Value = 9808.749899999999
$ws.Range("L107").Value = 4552.799999999999
$ws.Range("M107").Value = -7888.749899999999
$ws.Range("N107").Value = -8392.799999999999
$ws.Range("H136").Value = 8037.7812
$ws.Range("I136").Value = 4006.7407
$ws.Range("J136").Value = 9615.145
$ws.Range("K136").Value = 12020.2221
$ws.Range("L136").Value = 28845.435
$ws.Range("M136").Value = -9470.222099999999
$ws.Range("N136").Value = -33945.435
